$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("F").Delete()
